$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 score-summary header cells keep their "mtitleStyle" (s=4) formatting ---
# Rows 10-12 first column ("No.", "Marking", "Total") adopt that same style (s=4),
# matching the other header-row cells. Copy format from A9 (already s=4).
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Score summary numbers: handle float/re-graded input ---
$ws.Range("B10").Value = 14
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 11
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 56
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "53/112"

# --- Remove the third Student/Correct-Ans column pair (G:H) entirely ---
$ws.Range("G15:H21").Clear()

# --- Collapse the second Student-Ans column (D) into a real answer column,
#     reusing the "normalStyle" (s=5) / "incorrectStyle" (s=6) formats already
#     present in the sheet, then drop the now-unused Correct-Ans values (E19:E40) ---
$ws.Range("B10").Copy()
$ws.Range("D16:D17").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Value = "Option C"

$ws.Range("C10").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D18").Value = "Option B"

$ws.Range("D19:E40").Clear()

# --- Fill in the Student-Ans values for column A, rows 18-39, reusing the
#     "normalStyle" (s=5) / "incorrectStyle" (s=6) formats as appropriate ---
$ws.Range("B10").Copy()
$ws.Range("A18:A22").PasteSpecial(-4122)
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("A32:A33").PasteSpecial(-4122)
$ws.Range("A35:A37").PasteSpecial(-4122)
$ws.Range("A39").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C10").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A18").Value = "Option B"
$ws.Range("A19").Value = "Option C"
$ws.Range("A20").Value = "Option B"
$ws.Range("A21").Value = "Option C"
$ws.Range("A22").Value = "Option D"
$ws.Range("A25").Value = "Option B"
$ws.Range("A30").Value = "Option B"
$ws.Range("A32").Value = "Option C"
$ws.Range("A33").Value = "Option D"
$ws.Range("A34").Value = "Option A"
$ws.Range("A35").Value = "Option D"
$ws.Range("A36").Value = "Option A"
$ws.Range("A37").Value = "Option A"
$ws.Range("A39").Value = "Option D"
